# Commit message text change:
#   git commit -m "modif postform"
#   -> git commit -m "modif SAV & commentaires"
#
# The word "postform" is wrapped in <w:proofErr w:type="spellStart"/> /
# <w:proofErr w:type="spellEnd"/> (a stale spell-check marker). Per the
# target diff, the replacement text ("SAV & commentaires") must not be
# wrapped in those proofErr markers, so they need to disappear along with
# the old word - a plain Find/Replace would keep them wrapping the new
# run. We therefore rebuild the whole paragraph's OOXML explicitly and
# push it back in with Range.InsertXML (InsertXML replaces the full
# contents of the Range it is called on, so we call it on the whole
# paragraph Range to keep the surrounding markup/attributes intact).

$d = $word.ActiveDocument

$oldFragmentText = "postform"
$target = $null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$oldFragmentText*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing '$oldFragmentText'"
}

$newParagraphXml = '<w:p w14:paraId="334760CC" w14:textId="0A5B364E" w:rsidR="00FB703E" w:rsidRDefault="00FB703E" w:rsidP="00FB703E">' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> commit -m &quot;</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidR="00572814"><w:t>modif</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="00572814"><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>SAV &amp; commentaires</w:t></w:r>' +
    '<w:r><w:t>&quot;</w:t></w:r>' +
    '</w:p>'

$null = $target.Range.InsertXML($newParagraphXml)
